$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 351, pushing the existing rows 351-405 down to 353-407.
$ws.Rows("351:352").Insert()

# ---- New row 351 ----
$ws.Cells.Item(351, 1).Value = 10
$ws.Cells.Item(351, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(351, 3).Value = "La Araucanía"
$ws.Cells.Item(351, 4).Value = "02/27/2023"
$ws.Cells.Item(351, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(351, 5).Value = 9
$ws.Cells.Item(351, 6).Value = 100112001
$ws.Cells.Item(351, 7).Value = "Berenjena"
$ws.Cells.Item(351, 8).Value = "Sin especificar"
$ws.Cells.Item(351, 9).Value = "Primera"
$ws.Cells.Item(351, 10).Value = 170
$ws.Cells.Item(351, 11).Value = 14000
$ws.Cells.Item(351, 12).Value = 15000
$ws.Cells.Item(351, 13).Value = 14529
$ws.Cells.Item(351, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(351, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(351, 16).Value = 363
$ws.Cells.Item(351, 17).Value = 40
$ws.Cells.Item(351, 18).Value = "Hortaliza"

# ---- New row 352 ----
$ws.Cells.Item(352, 1).Value = 10
$ws.Cells.Item(352, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(352, 3).Value = "La Araucanía"
$ws.Cells.Item(352, 4).Value = "02/27/2023"
$ws.Cells.Item(352, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(352, 5).Value = 9
$ws.Cells.Item(352, 6).Value = 100112001
$ws.Cells.Item(352, 7).Value = "Berenjena"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Primera"
$ws.Cells.Item(352, 10).Value = 180
$ws.Cells.Item(352, 11).Value = 12000
$ws.Cells.Item(352, 12).Value = 15000
$ws.Cells.Item(352, 13).Value = 13333
$ws.Cells.Item(352, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(352, 15).Value = "Región del Maule"
$ws.Cells.Item(352, 16).Value = 333
$ws.Cells.Item(352, 17).Value = 40
$ws.Cells.Item(352, 18).Value = "Hortaliza"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
